$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: VALOR -> TOTAL, MODO -> TIPO
$ws.Range("B1").Value = "TOTAL"
$ws.Range("C1").Value = "TIPO"

# Update the active selection to C1 (as seen in saved file)
$ws.Range("C1").Select()
